# Weekly update: two new price records inserted at the top of the
# "Vega Modelo de Temuco - Mango" block (rows 651-681), pushing the
# existing rows down by two (now rows 653-683).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the current row 651 (twice,
# inserting at the same index pushes the previously-inserted row down).
$ws.Rows.Item(651).Insert()
$ws.Rows.Item(651).Insert()

# --- New row 651 --------------------------------------------------
$ws.Cells.Item(651, 1).Value = 10
$ws.Cells.Item(651, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(651, 3).Value = "La Araucanía"
$ws.Cells.Item(651, 4).Value = 45267
$ws.Cells.Item(651, 5).Value = 9
$ws.Cells.Item(651, 6).Value = "Fruta"
$ws.Cells.Item(651, 7).Value = 100108
$ws.Cells.Item(651, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(651, 9).Value = 100108002
$ws.Cells.Item(651, 10).Value = "Mango"
$ws.Cells.Item(651, 11).Value = "Sin especificar"
$ws.Cells.Item(651, 12).Value = "Primera"
$ws.Cells.Item(651, 13).Value = 305
$ws.Cells.Item(651, 14).Value = 12000
$ws.Cells.Item(651, 15).Value = 13000
$ws.Cells.Item(651, 16).Value = 12410
$ws.Cells.Item(651, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(651, 18).Value = "Perú"
$ws.Cells.Item(651, 19).Value = 3102
$ws.Cells.Item(651, 20).Value = 4

# --- New row 652 --------------------------------------------------
$ws.Cells.Item(652, 1).Value = 10
$ws.Cells.Item(652, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(652, 3).Value = "La Araucanía"
$ws.Cells.Item(652, 4).Value = 45267
$ws.Cells.Item(652, 5).Value = 9
$ws.Cells.Item(652, 6).Value = "Fruta"
$ws.Cells.Item(652, 7).Value = 100108
$ws.Cells.Item(652, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(652, 9).Value = 100108002
$ws.Cells.Item(652, 10).Value = "Mango"
$ws.Cells.Item(652, 11).Value = "Sin especificar"
$ws.Cells.Item(652, 12).Value = "Segunda"
$ws.Cells.Item(652, 13).Value = 250
$ws.Cells.Item(652, 14).Value = 10000
$ws.Cells.Item(652, 15).Value = 10000
$ws.Cells.Item(652, 16).Value = 10000
$ws.Cells.Item(652, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(652, 18).Value = "Perú"
$ws.Cells.Item(652, 19).Value = 2500
$ws.Cells.Item(652, 20).Value = 4
